$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.14907911348512926
$ws.Range("A2").Value = -0.046093148373724802
$ws.Range("A3").Value = -0.0089999995809524336
$ws.Range("A4").Value = 0.28399364616292999
$ws.Range("A5").Value = -0.0059999995917641158
$ws.Range("A6").Value = -0.0059999995741719658
$ws.Range("A7").Value = -0.01999999950045428
$ws.Range("A8").Value = -0.019999999496006282
$ws.Range("A9").Value = -0.0059999995623325475
$ws.Range("A10").Value = -0.0059999995574955278
$ws.Range("A11").Value = -0.0044999995648744573
$ws.Range("A12").Value = -0.0059999995561237363
$ws.Range("A13").Value = -0.0059999995516681892
$ws.Range("A14").Value = -0.011999999519665572
$ws.Range("A15").Value = 0.0093046583664024851
$ws.Range("A16").Value = -0.0059999995485335855
$ws.Range("A17").Value = -0.0059999995466260003
$ws.Range("A18").Value = -0.0089999995308014391
$ws.Range("A19").Value = -0.0089999995818179634
$ws.Range("A20").Value = 0.005685456248009757
$ws.Range("A21").Value = -0.062678654751868734
$ws.Range("A22").Value = -0.0089999995692084944
$ws.Range("A23").Value = -0.0089999995778500264
$ws.Range("A24").Value = -0.041999999401166832
$ws.Range("A25").Value = -0.041999999398070642
$ws.Range("A26").Value = -0.0059999995727473276
$ws.Range("A27").Value = -0.0059999995704944631
$ws.Range("A28").Value = -0.0059999995609594237
$ws.Range("A29").Value = -0.011999999523604643
$ws.Range("A30").Value = -0.019999999479504815
$ws.Range("A31").Value = -0.014999999500840744
$ws.Range("A32").Value = -0.020999999469349717
$ws.Range("A33").Value = -0.0059999995462440836
